# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Mon Apr 10 21:14:02 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    # Force text storage so numeric-looking strings (e.g. "6.250", "1.230")
    # keep their exact literal digits/trailing zeros instead of Excel
    # auto-coercing them into a Double and losing precision/formatting.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$ws.Range("D2").Value = "29.239.06"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "1.893.40"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  -0.21%  "
Set-TextCell $ws "D5" "314.54"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  -0.33%  "
Set-TextCell $ws "D7" "0.5142"
$ws.Range("E7").Value = "  +1.02%  "
Set-TextCell $ws "D8" "0.3912"
$ws.Range("E8").Value = "  -0.61%  "
Set-TextCell $ws "D9" "0.08408"
$ws.Range("E9").Value = "  +0.44%  "
Set-TextCell $ws "D10" "42.24"
$ws.Range("E10").Value = "  +1.73%  "
Set-TextCell $ws "D11" "1.113"
$ws.Range("E11").Value = "  +0.67%  "
Set-TextCell $ws "D12" "6.250"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.891.60"
$ws.Range("E13").Value = "  +0.96%  "
Set-TextCell $ws "D14" "20.62"
$ws.Range("E14").Value = "  +0.99%  "
Set-TextCell $ws "D15" "7.308"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("E16").Value = "  -0.21%  "
Set-TextCell $ws "D17" "93.08"
$ws.Range("E17").Value = "  +2.48%  "
Set-TextCell $ws "D18" "0.00001105"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("E19").Value = "  +0.45%  "
Set-TextCell $ws "D20" "17.82"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("E21").Value = "  -0.39%  "
Set-TextCell $ws "D22" "6.005"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "29.247.11"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("E24").Value = "  +0.26%  "
Set-TextCell $ws "D25" "2.218"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "2.107.52"
$ws.Range("E26").Value = "  +0.84%  "
Set-TextCell $ws "D27" "159.32"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +2.34%  "
Set-TextCell $ws "D30" "127.74"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D31" "0.1046"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D32" "1.056"
$ws.Range("E32").Value = "  +1.08%  "
Set-TextCell $ws "D33" "6.174"
$ws.Range("E33").Value = "  +7.14%  "
Set-TextCell $ws "D34" "3.661"
$ws.Range("E34").Value = "  +2.03%  "
Set-TextCell $ws "D35" "0.02479"
$ws.Range("E35").Value = "  +1.98%  "
Set-TextCell $ws "D36" "0.06565"
$ws.Range("E36").Value = "  +1.44%  "
Set-TextCell $ws "D37" "9.028"
$ws.Range("E37").Value = "  +1.78%  "
Set-TextCell $ws "D38" "0.2188"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws "D39" "5.179"
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D40" "1.226"
$ws.Range("E40").Value = "  +3.11%  "
Set-TextCell $ws "D41" "0.6493"
$ws.Range("E41").Value = "  +1.20%  "
Set-TextCell $ws "D42" "1.230"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("E43").Value = "  +1.08%  "
Set-TextCell $ws "D44" "0.6043"
$ws.Range("E44").Value = "  +0.08%  "
Set-TextCell $ws "D45" "13.19"
$ws.Range("E45").Value = "  +1.08%  "
Set-TextCell $ws "D46" "3.670"
$ws.Range("E46").Value = "  -0.63%  "
Set-TextCell $ws "D47" "2.048"
$ws.Range("E47").Value = "  +2.15%  "
Set-TextCell $ws "D48" "1.228"
$ws.Range("E48").Value = "  +2.17%  "
Set-TextCell $ws "D49" "123.41"
$ws.Range("E49").Value = "  +1.35%  "
Set-TextCell $ws "D50" "1.173"
$ws.Range("E50").Value = "  -1.22%  "
Set-TextCell $ws "D51" "77.35"
$ws.Range("E51").Value = "  +0.80%  "
